$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Duplicate the S1 notes table (rows 2-39) as an S2 notes table (rows 40-77):
# same Libelle/Matiere/Coefficient/Fichier, but Dossier_Semestre becomes
# "Notes_S2" and the Unite_d_Enseignement "UE1.x" becomes "UE2.x".
for ($i = 0; $i -lt 38; $i++) {
    $srcRow = 2 + $i
    $dstRow = 40 + $i
    $src = $ws.Range("A" + $srcRow + ":G" + $srcRow)
    $dst = $ws.Range("A" + $dstRow + ":G" + $dstRow)
    $src.Copy($dst)

    $ws.Range("B" + $dstRow).Value = "Notes_S2"

    $oldUE = $ws.Range("E" + $dstRow).Value()
    $newUE = $oldUE -replace "UE1\.", "UE2."
    $ws.Range("E" + $dstRow).Value = $newUE
}

# Match the author's final view state: zoomed to 120%, scrolled down so
# row 28 is at the top, with the newly-added E40:E77 column selected.
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.Zoom = 120
$ws.Range("E40:E77").Select()
